$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.514.95"
$ws.Range("D3").Value = "1.977.32"
$ws.Range("E3").Value = "  -3.84%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("E6").Value = "  -4.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "56.78"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.79%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "58.29"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.26%  "
$ws.Range("E10").Value = "  -0.70%  "
$ws.Range("E11").Value = "  -2.49%  "
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("D15").Value = "2.264.32"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("E16").Value = "  -2.76%  "
$ws.Range("D17").Value = "1.980.24"
$ws.Range("E17").Value = "  -4.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.70%  "
$ws.Range("D19").Value = "35.466.10"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  -2.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "232.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  +17.73%  "
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.09%  "
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("E31").Value = "  -4.95%  "
$ws.Range("E32").Value = "  -10.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0953"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +16.23%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  +8.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("E38").Value = "  -4.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  -2.65%  "
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "91.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.90%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0879"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.42%  "
$ws.Range("D48").Value = "1.365.49"
$ws.Range("E48").Value = "  -3.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.73%  "
